$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# --- Cells changing from text-placeholder to numeric: copy style from an untouched numeric cell, then set value ---
$ws.Range("G15").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("G15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("N15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 100

# --- Cells changing from numeric to text-placeholder "0": copy whole cell (style+value) from untouched placeholder cell C22 ---
$ws.Range("C22").Copy($ws.Range("D25"))
$ws.Range("C22").Copy($ws.Range("C28"))
$ws.Range("C22").Copy($ws.Range("D29"))
$ws.Range("C22").Copy($ws.Range("D30"))

# --- Cells changing from numeric to text-placeholder "***.*": copy whole cell (style+value) from untouched placeholder cell N23 ---
$ws.Range("N23").Copy($ws.Range("E25"))
$ws.Range("N23").Copy($ws.Range("E29"))
$ws.Range("N23").Copy($ws.Range("E30"))

# --- Simple numeric updates (style unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 50
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 23
$ws.Range("H16").Value = 64.285714285714
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 31.818181818181
$ws.Range("L16").Value = 52.631578947368
$ws.Range("M16").Value = 11.538461538461
$ws.Range("N16").Value = -64.197530864197
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 133.333333333333
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -35.483870967741
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 43
$ws.Range("K17").Value = -11.627906976744
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 153.333333333333
$ws.Range("N17").Value = -39.682539682539
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = -37.5
$ws.Range("M18").Value = -9.090909090909
$ws.Range("N18").Value = -74.358974358974
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 7.692307692307
$ws.Range("I19").Value = 38
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = -2.564102564102
$ws.Range("L19").Value = -30.90909090909
$ws.Range("M19").Value = 58.333333333333
$ws.Range("N19").Value = -50
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = -60
$ws.Range("L20").Value = -83.333333333333
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -95.348837209302
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -3.529411764705
$ws.Range("I21").Value = 121
$ws.Range("J21").Value = 125
$ws.Range("K21").Value = -3.2
$ws.Range("L21").Value = -16.551724137931
$ws.Range("M21").Value = 49.382716049382
$ws.Range("N21").Value = -61.708860759493
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = -15.78947368421
$ws.Range("L23").Value = -27.272727272727
$ws.Range("M23").Value = 88.235294117647
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 66.666666666666
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 39.682539682539
$ws.Range("I24").Value = 107
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 33.75
$ws.Range("L24").Value = 35.443037974683
$ws.Range("M24").Value = 50.704225352112
$ws.Range("C25").Value = 8
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 212.5
$ws.Range("I25").Value = 29
$ws.Range("J25").Value = 11
$ws.Range("K25").Value = 163.636363636364
$ws.Range("L25").Value = 123.076923076923
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -45.454545454545
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -10.25641025641
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 55
$ws.Range("K26").Value = -9.090909090909
$ws.Range("L26").Value = -15.254237288135
$ws.Range("M26").Value = -21.875
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 2
$ws.Range("L27").Value = -42.857142857142
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = -16.666666666666
